$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "STUDY NAME"
$ws.Range("C1").Value = "YEAR"

# Study names (column B)
$ws.Range("B2").Value = "Promoting Walking and Bicycling: Assessing the Evidence to Assist Planners"
$ws.Range("D1").Value = "Study Type"
$ws.Range("B3").Value = "Shaping the City for Walking and Cycling: A Case Study of Lancaster"

# Study type (column D)
$ws.Range("D3").Value = "Case Study"
$ws.Range("D4").Value = "Case Study"
$ws.Range("D2").Value = "Studies Review"
$ws.Range("B4").Value = "The Impact of Local Transport Systems on Green Infrastructure – Policy Versus Reality: The Case of Poznan, Poland"

# ID and Year columns
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 2010
$ws.Range("A3").Value = 2
$ws.Range("C3").Value = 2010
$ws.Range("A4").Value = 3
$ws.Range("C4").Value = 2015

# Apply Times (black) font to the Study Name cells
foreach ($addr in @("B2", "B3", "B4")) {
    $cell = $ws.Range($addr)
    $cell.Font.Color = 0
    $cell.Font.Name = "Times"
}

# Column widths (best-fit equivalents)
$ws.Columns.Item(2).ColumnWidth = 97.16666666666667
$ws.Columns.Item(4).ColumnWidth = 12.830729166666666

# Final selection, as left by the editing session
$ws.Range("A5").Select()
